$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 16 (before SPEC 01), shifting existing rows down
$ws.Rows.Item(16).Resize(2).Insert()

# Set the new values
$ws.Range("A16").Value = "ACCOUNTING"
$ws.Range("A17").Value = "BILLING"

# Update selection to A18
$ws.Range("A18").Select()
